# Updates cryptocurrency price (column D) and volume change (column E) values
# to reflect the latest scrape, as captured in the commit's OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D sometimes holds price text that is also a valid number
# (e.g. "0.9999", "84.40"). Writing such text via .Value would make Excel
# silently reinterpret it as a numeric value and normalize/round it
# (".40" -> ".4", trailing zeros dropped, etc.), which would not match the
# original inline-string cell content. To avoid this, those specific cells
# are temporarily switched to a text ("@") number format before the value
# is written, then restored to the default "Normal" style so no visible
# formatting change is introduced.
$textForcedCells = @(
    "D4", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D15", "D16", "D17", "D20", "D22", "D23", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# New values for every changed cell (Price in column D, Volume(1h) in column E)
$newValues = @{
    'D2' = '29.907.65'
    'E2' = '  +2.81%  '
    'D3' = '1.861.44'
    'E3' = '  +2.27%  '
    'D4' = '0.9995'
    'E4' = '  +0.14%  '
    'D5' = '246.47'
    'E5' = '  +2.15%  '
    'D6' = '0.6362'
    'E6' = '  +3.36%  '
    'D8' = '0.3008'
    'E8' = '  +4.26%  '
    'D9' = '0.07495'
    'E9' = '  +2.33%  '
    'D10' = '24.71'
    'E10' = '  +8.25%  '
    'D11' = '0.07683'
    'E11' = '  +0.35%  '
    'D12' = '1.867.07'
    'E12' = '  +2.47%  '
    'D13' = '5.054'
    'E13' = '  +2.30%  '
    'E14' = '  +4.72%  '
    'D15' = '84.40'
    'E15' = '  +3.59%  '
    'D16' = '0.000009392'
    'E16' = '  +4.78%  '
    'D17' = '6.117'
    'E17' = '  +4.90%  '
    'D18' = '29.849.05'
    'D19' = '2.107.44'
    'E19' = '  +2.15%  '
    'D20' = '240.03'
    'E20' = '  +1.12%  '
    'E21' = '  +2.10%  '
    'D22' = '0.9999'
    'E22' = '  -0.05%  '
    'D23' = '7.365'
    'E23' = '  +3.74%  '
    'E24' = '  +0.12%  '
    'D25' = '159.26'
    'E25' = '  +1.01%  '
    'D26' = '0.1422'
    'E26' = '  +1.34%  '
    'D27' = '8.585'
    'E27' = '  +1.91%  '
    'E28' = '  +2.44%  '
    'D29' = '1.507'
    'E29' = '  +2.17%  '
    'D30' = '0.06051'
    'D31' = '1.270'
    'E31' = '  +5.30%  '
    'D32' = '4.153'
    'E32' = '  +1.68%  '
    'D33' = '4.146'
    'E33' = '  +1.36%  '
    'D34' = '1.875'
    'E34' = '  +3.65%  '
    'D35' = '1.164'
    'E35' = '  +3.14%  '
    'D36' = '0.7303'
    'E36' = '  -0.36%  '
    'D37' = '2.614'
    'E37' = '  -0.08%  '
    'E38' = '  +1.45%  '
    'D39' = '1.223.11'
    'E39' = '  +1.24%  '
    'D40' = '0.01789'
    'E40' = '  +1.76%  '
    'D41' = '6.312'
    'E41' = '  -0.96%  '
    'D42' = '0.9237'
    'E42' = '  +3.52%  '
    'D43' = '1.001'
    'E43' = '  -0.01%  '
    'D44' = '2.018.48'
    'E44' = '  +2.51%  '
    'D45' = '102.37'
    'E45' = '  +1.60%  '
    'D46' = '66.56'
    'E46' = '  +3.34%  '
    'D47' = '0.00000000122'
    'E47' = '  +1.19%  '
    'D48' = '0.5085'
    'E48' = '  +0.16%  '
    'D49' = '9.358'
    'E49' = '  +3.97%  '
    'D50' = '0.4096'
    'E50' = '  +2.78%  '
    'D51' = '0.1144'
    'E51' = '  +3.02%  '
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}

# Restore the default cell style on the cells we temporarily reformatted,
# so the only lasting change is the cell content, matching the source diff.
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).Style = "Normal"
}
